# "Update countries & provincias Spain"
# Refresh the COVID-19 country snapshot on sheet "Pais":
#  - a handful of countries swap rank/row position because their
#    "Casos totales" (column B) changed relative to their neighbours
#    (e.g. Irak now outranks Irlanda & Republica Dominicana, etc.)
#  - most rows simply get refreshed totals for columns B..H
#  - the "last updated" timestamp banner in A1 is bumped

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each tuple is (row, column, new value) using 1-based column numbers
# (A=1 ... H=8). Only cells whose value actually changes are listed.
$updates = @(
    @(1, 1, "Datos actualizados a 18 de Junio de 2020 a las 19:12"),
    @(4, 2, 2246676),
    @(4, 3, 12205),
    @(4, 4, 920427),
    @(4, 5, 1206019),
    @(4, 7, 289),
    @(4, 8, 120230),
    @(5, 2, 965512),
    @(5, 3, 5203),
    @(5, 5, 415163),
    @(5, 7, 177),
    @(5, 8, 46842),
    @(7, 2, 377122),
    @(7, 3, 9858),
    @(7, 4, 200358),
    @(7, 5, 164263),
    @(7, 7, 239),
    @(7, 8, 12501),
    @(9, 2, 292348),
    @(9, 3, 585),
    @(11, 2, 238159),
    @(11, 3, 331),
    @(11, 4, 180544),
    @(11, 5, 23101),
    @(11, 7, 66),
    @(11, 8, 34514),
    @(12, 2, 225103),
    @(12, 3, 4475),
    @(12, 5, 65030),
    @(12, 7, 226),
    @(12, 8, 3841),
    @(14, 2, 189565),
    @(14, 3, 61),
    @(14, 4, 174100),
    @(14, 5, 6532),
    @(14, 7, 6),
    @(14, 8, 8933),
    @(21, 2, 100146),
    @(21, 3, 293),
    @(21, 4, 62445),
    @(21, 5, 29402),
    @(21, 7, 45),
    @(21, 8, 8299),
    @(31, 2, 49097),
    @(31, 3, 607),
    @(31, 4, 24123),
    @(31, 5, 20887),
    @(31, 7, 80),
    @(31, 8, 4087),
    @(44, 1, "Irak"),
    @(44, 2, 25717),
    @(44, 3, 1463),
    @(44, 4, 11333),
    @(44, 5, 13528),
    @(44, 7, 83),
    @(44, 8, 856),
    @(45, 1, "Irlanda"),
    @(45, 2, 25341),
    @(45, 3, 0),
    @(45, 4, 22698),
    @(45, 5, 933),
    @(45, 7, 0),
    @(45, 8, 1710),
    @(46, 1, "Republica Dominicana"),
    @(46, 2, 24645),
    @(46, 3, 540),
    @(46, 4, 14293),
    @(46, 5, 9717),
    @(46, 7, 2),
    @(46, 8, 635),
    @(62, 1, "Argelia"),
    @(62, 2, 11385),
    @(62, 3, 117),
    @(62, 4, 8078),
    @(62, 5, 2496),
    @(62, 7, 12),
    @(62, 8, 811),
    @(63, 1, "Azerbaiyan"),
    @(63, 2, 11329),
    @(63, 3, 338),
    @(63, 4, 6192),
    @(63, 5, 4998),
    @(63, 7, 6),
    @(63, 8, 139),
    @(66, 2, 10230),
    @(66, 3, 68),
    @(66, 4, 7436),
    @(66, 5, 2460),
    @(66, 7, 1),
    @(66, 8, 334),
    @(83, 2, 4557),
    @(83, 3, 12),
    @(83, 4, 3527),
    @(83, 5, 987),
    @(121, 2, 1330),
    @(121, 3, 22),
    @(121, 4, 717),
    @(121, 5, 600),
    @(122, 2, 1272),
    @(122, 3, 23),
    @(122, 4, 710),
    @(122, 5, 511),
    @(128, 4, 818),
    @(128, 5, 148),
    @(128, 7, 1),
    @(128, 8, 19),
    @(145, 1, "Estado de Palestina"),
    @(145, 2, 599),
    @(145, 3, 44),
    @(145, 4, 415),
    @(145, 5, 181),
    @(145, 7, 0),
    @(145, 8, 3),
    @(146, 1, "Benin"),
    @(146, 2, 597),
    @(146, 3, 25),
    @(146, 4, 238),
    @(146, 5, 348),
    @(146, 7, 2),
    @(146, 8, 11),
    @(157, 1, "Montenegro"),
    @(157, 3, 4),
    @(157, 4, 315),
    @(157, 5, 13),
    @(157, 8, 9),
    @(158, 1, "Mauricio"),
    @(158, 2, 337),
    @(158, 4, 325),
    @(158, 5, 2),
    @(158, 8, 10),
    @(159, 1, "Isla de Man"),
    @(159, 2, 336),
    @(159, 4, 312),
    @(159, 5, 0),
    @(159, 8, 24),
    @(160, 4, 187),
    @(160, 5, 70),
    @(163, 1, "Comoras"),
    @(163, 2, 210),
    @(163, 3, 13),
    @(163, 4, 129),
    @(163, 5, 76),
    @(163, 7, 2),
    @(163, 8, 5),
    @(164, 1, "Mongolia"),
    @(164, 2, 201),
    @(164, 3, 4),
    @(164, 5, 74),
    @(164, 8, 0),
    @(202, 1, "Fiyi"),
    @(203, 1, "Dominica"),
    @(206, 1, "Groenlandia"),
    @(207, 1, "Islas Malvinas"),
    @(208, 1, "Islas Turcas y Caicos"),
    @(208, 4, 11),
    @(208, 8, 1),
    @(209, 1, "Santa Sede"),
    @(209, 4, 12),
    @(209, 8, 0),
)

foreach ($u in $updates) {
    $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
}